# Update the division-answer table: replace each cell's text with its
# new value. Several old values repeat (e.g. "64÷6=10, 4" appears twice)
# so we address each cell positionally by (row, column) instead of using
# a global Find/Replace, which would not be able to distinguish the
# duplicate occurrences or target the correct one.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# (row, col, expectedOldText, newText)
$updates = @(
    @(1, 1, "21÷4=5, 1",  "24÷5=4, 4"),
    @(1, 2, "43÷5=8, 3",  "76÷5=15, 1"),
    @(1, 3, "79÷8=9, 7",  "56÷6=9, 2"),
    @(1, 4, "64÷6=10, 4", "64÷2=32, 0"),
    @(1, 5, "17÷9=1, 8",  "75÷5=15, 0"),

    @(5, 1, "29÷4=7, 1",  "75÷3=25, 0"),
    @(5, 2, "93÷9=10, 3", "61÷6=10, 1"),
    @(5, 3, "27÷8=3, 3",  "99÷4=24, 3"),
    @(5, 4, "64÷6=10, 4", "42÷6=7, 0"),
    @(5, 5, "29÷2=14, 1", "21÷4=5, 1"),

    @(9, 1, "90÷7=12, 6", "27÷2=13, 1"),
    @(9, 2, "59÷8=7, 3",  "22÷4=5, 2"),
    @(9, 3, "50÷2=25, 0", "11÷2=5, 1"),
    @(9, 4, "43÷3=14, 1", "70÷8=8, 6"),
    @(9, 5, "49÷6=8, 1",  "48÷8=6, 0"),

    @(13, 1, "92÷8=11, 4", "44÷4=11, 0"),
    @(13, 2, "66÷8=8, 2",  "55÷9=6, 1"),
    @(13, 3, "83÷9=9, 2",  "77÷5=15, 2"),
    @(13, 4, "58÷4=14, 2", "20÷8=2, 4"),
    @(13, 5, "21÷6=3, 3",  "36÷6=6, 0"),

    @(17, 1, "86÷6=14, 2", "82÷5=16, 2"),
    @(17, 2, "83÷8=10, 3", "41÷6=6, 5"),
    @(17, 3, "58÷7=8, 2",  "20÷7=2, 6"),
    @(17, 4, "92÷9=10, 2", "50÷4=12, 2"),
    @(17, 5, "70÷5=14, 0", "23÷3=7, 2")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $oldText = $u[2]
    $newText = $u[3]

    # Re-fetch the table/cell fresh for every single update: after a
    # text-length-changing edit, previously obtained Cell/Range handles
    # can go stale, so always address cells via a newly retrieved
    # collection right before use.
    $cell = $d.Tables.Item(1).Cell($row, $col)

    # Sanity-check we are about to overwrite the expected current value
    # (cell text includes a trailing end-of-cell mark, so use -like).
    $current = $cell.Range.Text
    if ($current -notlike "$oldText*") {
        throw "Unexpected text in row $row col $col`: [$current] (expected to start with [$oldText])"
    }

    # Assigning directly to Range.Text replaces only the run's text
    # while the existing run/paragraph formatting (font, size, etc.)
    # is preserved.
    $cell.Range.Text = $newText
}
